$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.210793852806091
$ws.Range("B1").Value = 1.470886588096619
$ws.Range("C1").Value = 1.832079887390137
$ws.Range("D1").Value = 1.672866702079773
$ws.Range("E1").Value = 1.550372362136841
